$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10; this shifts the existing rows 10-22 down to
# 11-23 (and the sheet's used range grows from A1:R22 to A1:R23), matching the
# diff which re-numbers every record from row 10 onward by one and appends a
# new row 23 containing what used to be row 22's data.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new weekly record.
$ws.Cells.Item(10, 1).Value = 5
$ws.Cells.Item(10, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(10, 3).Value = "Maule"
$ws.Cells.Item(10, 4).Value = 44614
$ws.Cells.Item(10, 5).Value = 7
$ws.Cells.Item(10, 6).Value = 100112043
$ws.Cells.Item(10, 7).Value = "Pepino dulce"
$ws.Cells.Item(10, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 300
$ws.Cells.Item(10, 11).Value = 15000
$ws.Cells.Item(10, 12).Value = 15000
$ws.Cells.Item(10, 13).Value = 15000
$ws.Cells.Item(10, 14).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(10, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(10, 16).Value = 833
$ws.Cells.Item(10, 17).Value = 18
$ws.Cells.Item(10, 18).Value = "Hortaliza"
